$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the generic "1" to the municipality name
$ws.Name = "საჩხერე"

# Remove the now spare blank spacer row (old row 3) entirely, shifting the
# rows below it up by one.
$ws.Range("A3").EntireRow.Delete()

# The row that used to hold "(მოსახლეობის აღწერის შედეგებით)" is kept but
# emptied out completely.
$ws.Range("A2:B2").Clear()

# Drop the 1989 and 2002 columns, keeping only the 2014 figures which slide
# into column B.
$ws.Range("B1:C1").EntireColumn.Delete()

# Match the saved selection/active cell from the source file.
$ws.Range("A2").Select()
